$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell value while forcing it to remain plain text
# (prevents Excel auto-converting numeric-looking strings like "1.00" into numbers)
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "66.338.36"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "3.079.25"
$ws.Range("E3").Value = "  -1.10%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "575.83"
$ws.Range("E5").Value = "  -0.33%  "
Set-TextValue $ws.Range("D6") "169.99"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.075.24"
$ws.Range("E8").Value = "  -1.11%  "
Set-TextValue $ws.Range("D9") "0.509"
$ws.Range("E9").Value = "  -2.19%  "
Set-TextValue $ws.Range("D10") "6.41"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  -2.81%  "
Set-TextValue $ws.Range("D13") "0.0000239"
$ws.Range("E13").Value = "  -2.53%  "
Set-TextValue $ws.Range("D14") "35.79"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "3.591.86"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "66.294.02"
$ws.Range("E17").Value = "  -0.75%  "
Set-TextValue $ws.Range("D18") "6.97"
$ws.Range("E18").Value = "  -2.86%  "
Set-TextValue $ws.Range("D19") "16.89"
$ws.Range("E19").Value = "  +3.58%  "
$ws.Range("D20").Value = "3.075.91"
$ws.Range("E20").Value = "  -1.14%  "
Set-TextValue $ws.Range("D21") "488.36"
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("E23").Value = "  -3.58%  "
Set-TextValue $ws.Range("D24") "82.79"
$ws.Range("E24").Value = "  -1.49%  "
Set-TextValue $ws.Range("D25") "12.71"
$ws.Range("E25").Value = "  -4.61%  "
$ws.Range("E26").Value = "  -2.90%  "
Set-TextValue $ws.Range("D27") "10.16"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  +0.06%  "
Set-TextValue $ws.Range("D29") "7.84"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("E31").Value = "  -2.01%  "
Set-TextValue $ws.Range("D32") "27.70"
$ws.Range("E32").Value = "  -2.95%  "
Set-TextValue $ws.Range("D33") "0.111"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("D34").Value = "0.0₃0911"
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("B36").Value = "Arweave"
$ws.Range("C36").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D36") "47.46"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D37") "5.60"
$ws.Range("E37").Value = "  -4.32%  "
Set-TextValue $ws.Range("D38") "0.949"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("D43").Value = "2.792.32"
$ws.Range("E43").Value = "  -1.23%  "
Set-TextValue $ws.Range("D44") "2.52"
$ws.Range("E44").Value = "  -1.16%  "
Set-TextValue $ws.Range("D45") "0.0347"
$ws.Range("E45").Value = "  -2.30%  "
Set-TextValue $ws.Range("D46") "366.86"
$ws.Range("E46").Value = "  -4.38%  "
Set-TextValue $ws.Range("D47") "134.56"
$ws.Range("E47").Value = "  -0.60%  "
Set-TextValue $ws.Range("D49") "24.45"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("E51").Value = "  -2.09%  "
